$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2..6), matching the diff's final state.
$data = @(
    @{ Row=2; D=44307; M=250; N=19000; O=20000; P=19500; Q="$/bandeja 18 kilos"; S=1083; T=18 },
    @{ Row=3; D=44323; M=270; N=21000; O=22000; P=21500; Q="$/bandeja 18 kilos"; S=1194; T=18 },
    @{ Row=4; D=44418; M=240; N=10000; O=11000; P=10500; Q="$/bandeja 10 kilos"; S=1050; T=10 },
    @{ Row=5; D=44291; M=200; N=17000; O=18000; P=17500; Q="$/bandeja 18 kilos"; S=972;  T=18 },
    @{ Row=6; D=44263; M=250; N=21000; O=22000; P=21500; Q="$/caja 18 kilos";    S=1194; T=18 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D    # D: Fecha
    $ws.Cells.Item($r, 13).Value = $item.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $item.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $item.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $item.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $item.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 19).Value = $item.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $item.T   # T: Kg / unidad
}
